$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (matching source data type)

# Row 2
$ws.Range("D2").Value = "62.451.53"
$ws.Range("E2").Value = "  -4.81%  "

# Row 3
$ws.Range("D3").Value = "3.214.74"
$ws.Range("E3").Value = "  -5.54%  "

# Row 4
$ws.Range("E4").Value = "  +0.51%  "

# Row 5
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "516.18"
$ws.Range("E5").Value = "  -3.22%  "

# Row 6
$ws.Range("B6").Value = "Solana"
$ws.Range("C6").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.81"
$ws.Range("E6").Value = "  -6.02%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.589"
$ws.Range("E7").Value = "  -3.57%  "

# Row 8
$ws.Range("E8").Value = "  +0.24%  "

# Row 9
$ws.Range("D9").Value = "3.210.15"
$ws.Range("E9").Value = "  -5.58%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.596"
$ws.Range("E10").Value = "  -5.49%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.13"
$ws.Range("E11").Value = "  -10.86%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.130"
$ws.Range("E12").Value = "  -3.98%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000250"
$ws.Range("E13").Value = "  -3.30%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.86"
$ws.Range("E14").Value = "  -5.57%  "

# Row 15
$ws.Range("D15").Value = "3.771.84"
$ws.Range("E15").Value = "  -3.99%  "

# Row 16
$ws.Range("D16").Value = "3.253.11"
$ws.Range("E16").Value = "  -3.72%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.114"
$ws.Range("E17").Value = "  -7.68%  "

# Row 18
$ws.Range("D18").Value = "62.714.05"
$ws.Range("E18").Value = "  -3.72%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.03"
$ws.Range("E19").Value = "  -3.61%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.89"
$ws.Range("E20").Value = "  -3.30%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.952"
$ws.Range("E21").Value = "  -3.00%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "362.83"
$ws.Range("E22").Value = "  -3.72%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.70"
$ws.Range("E23").Value = "  -1.48%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.80"
$ws.Range("E24").Value = "  -4.02%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.89"
$ws.Range("E25").Value = "  -1.17%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.87"
$ws.Range("E26").Value = "  +5.64%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.13"
$ws.Range("E27").Value = "  +5.38%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.59"
$ws.Range("E28").Value = "  -3.49%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.09"
$ws.Range("E29").Value = "  -5.12%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.10"
$ws.Range("E30").Value = "  -5.40%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "650.14"
$ws.Range("E31").Value = "  -5.20%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.05"
$ws.Range("E32").Value = "  -6.29%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.35"
$ws.Range("E33").Value = "  -7.14%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.02"
$ws.Range("E34").Value = "  -2.01%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.103"
$ws.Range("E35").Value = "  -2.88%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.18"
$ws.Range("E36").Value = "  -7.34%  "

# Row 37
$ws.Range("E37").Value = "  -0.08%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.12"
$ws.Range("E38").Value = "  -2.46%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.371"
$ws.Range("E39").Value = "  -5.29%  "

# Row 40
$ws.Range("E40").Value = "  +0.91%  "

# Row 41
$ws.Range("D41").Value = "0.0₃0703"
$ws.Range("E41").Value = "  +10.67%  "

# Row 42
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "2.870.51"
$ws.Range("E42").Value = "  -0.31%  "

# Row 43
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.121"
$ws.Range("E43").Value = "  -6.23%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.47"
$ws.Range("E44").Value = "  +4.15%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.62"
$ws.Range("E45").Value = "  -1.12%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0385"
$ws.Range("E46").Value = "  -2.74%  "

# Row 47
$ws.Range("B47").Value = "ThetaToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.53"
$ws.Range("E47").Value = "  -9.38%  "

# Row 48
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.76"
$ws.Range("E48").Value = "  +6.10%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.96"
$ws.Range("E49").Value = "  +3.35%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "134.67"
$ws.Range("E50").Value = "  -0.97%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.122"
$ws.Range("E51").Value = "  -3.56%  "
